# Orders sheet: row 9 ("Ipoh Coff") was a near-duplicate / mistyped product
# name; fix it to the canonical "Ipoh Coffee" product already used elsewhere
# in the sheet (e.g. A3). Once no cell references the old "Ipoh Coff" shared
# string, it drops out of sharedStrings.xml and every later shared-string
# index shifts down by one - which is what the rest of the diff reflects.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$ws.Range("A9").Value = "Ipoh Coffee"

# Leave the selection on A9, matching the saved cursor position.
$ws.Range("A9").Select()
